$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this pushes the existing rows 6..82
# down to 7..83, matching the rest of the dataset unchanged.
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new weekly data point.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44756
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112040
$ws.Cells.Item(6, 7).Value = "Cilantro"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 300
$ws.Cells.Item(6, 11).Value = 3500
$ws.Cells.Item(6, 12).Value = 4000
$ws.Cells.Item(6, 13).Value = 3750
$ws.Cells.Item(6, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 1875
$ws.Cells.Item(6, 17).Value = 2
$ws.Cells.Item(6, 18).Value = "Hortaliza"
